$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C10: 18 -> 1
$ws.Range("C10").Value = 1

# B11: fix mistaken "1" label back to "R40" (matches R10/R20/R30 pattern)
$ws.Range("B11").Value = "R40"
